# Agrego CP03 y CP07
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content changes ---
# CP03 LoginUserInvalido -> CP03LoginUserInvalido, and add the error message text
$ws.Range("A5").Value = "CP03LoginUserInvalido"
$ws.Range("D5").Value = "Login was unsuccessful. Please correct the errors and try again."

# CP07 SubNewsletter -> CP07SubNewsletter, add test data and confirmation text
$ws.Range("A9").Value = "CP07SubNewsletter"
$ws.Range("C9").Value = 123456
$ws.Range("D9").Value = "Thank you for signing up! A verification email has been sent. We appreciate your interest."

# New hyperlink on B9 (mailto link, like the other test-data rows)
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:prueba_1@gmail.com")

# --- Formatting changes ---
# Header row -> blue fill
$ws.Range("A1:F1").Interior.Color = 15773696

# Rows 2-5 and the new row 9 (CP00-CP03 and CP07) -> green fill
$ws.Range("A2:F5").Interior.Color = 5296274
$ws.Range("A9:F9").Interior.Color = 5296274

# Widen column D to fit the new long text
$ws.Columns("D").ColumnWidth = 81.7

# --- View state ---
$ws.Range("B14").Select()
